$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "https://www.impact.science/"
$ws.Range("B2").Value = "_gat_gtag_UA_140842239_1"
$ws.Range("F2").Value = "2021-08-16 21:21:01 India Standard Time"
$ws.Range("H2").Value = "2023-08-16 21:20:01 India Standard Time"
$ws.Range("J2").Value = 1692201002

# Row 3 updates
$ws.Range("A3").Value = "https://lifesciences.cactusglobal.com/"
$ws.Range("B3").Value = "_ga_MNGCCS5STP"
$ws.Range("F3").Value = "2023-08-16 21:20:23 India Standard Time"
$ws.Range("G3").Value = "Yes"
$ws.Range("H3").Value = "2023-08-16 21:20:24 India Standard Time"
$ws.Range("I3").Value = "Yes"
$ws.Range("J3").Value = 1692201025

# Update active selection
$ws.Range("G6").Select()
